$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# Row 10 and Row 12 swap their "entry" data (everything except the
# ID/DATE in columns A/B, the ESPORTS category in column G, and the
# running-balance formulas in D/F/M which recompute automatically).
# Capture the original values first (use Value2 - Value's getter is not
# reliable for reads through this bridge).
$c10 = $ws.Range("C10").Value2
$e10 = $ws.Range("E10").Value2
$h10 = $ws.Range("H10").Value2
$i10 = $ws.Range("I10").Value2
$j10 = $ws.Range("J10").Value2
$k10 = $ws.Range("K10").Value2
$l10 = $ws.Range("L10").Value2

$c12 = $ws.Range("C12").Value2
$e12 = $ws.Range("E12").Value2
$h12 = $ws.Range("H12").Value2
$i12 = $ws.Range("I12").Value2
$j12 = $ws.Range("J12").Value2
$k12 = $ws.Range("K12").Value2
$l12 = $ws.Range("L12").Value2

# Write row 12's original entry into row 10
$ws.Range("C10").Value2 = $c12
$ws.Range("E10").Value2 = $e12
$ws.Range("H10").Value2 = $h12
$ws.Range("I10").Value2 = $i12
$ws.Range("J10").Value2 = $j12
$ws.Range("K10").Value2 = $k12
$ws.Range("L10").Value2 = $l12

# Write row 10's original entry into row 12
$ws.Range("C12").Value2 = $c10
$ws.Range("E12").Value2 = $e10
$ws.Range("H12").Value2 = $h10
$ws.Range("I12").Value2 = $i10
$ws.Range("J12").Value2 = $j10
$ws.Range("K12").Value2 = $k10
$ws.Range("L12").Value2 = $l10

# Move the active cell selection from E8 to H6
$ws.Range("H6").Select()
